$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 159, pushing the existing rows 159-168 down
# to 161-170 (Excel copies formatting, including the date-style D column,
# from the row above on insert - matches the style used throughout column D).
$ws.Rows("159:160").Insert()

# New row 159: Terminal Hortofrutícola Agro Chillán, Ñuble, Mandarina - Clemenuless, Primera
$ws.Range("A159").Value = 7
$ws.Range("B159").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C159").Value = "Ñuble"
$ws.Range("D159").Value = 44753
$ws.Range("E159").Value = 16
$ws.Range("F159").Value = "Fruta"
$ws.Range("G159").Value = 100102
$ws.Range("H159").Value = "Cítricos"
$ws.Range("I159").Value = 100102004
$ws.Range("J159").Value = "Mandarina"
$ws.Range("K159").Value = "Clemenuless"
$ws.Range("L159").Value = "Primera"
$ws.Range("M159").Value = 160
$ws.Range("N159").Value = 7000
$ws.Range("O159").Value = 7500
$ws.Range("P159").Value = 7250
$ws.Range("Q159").Value = "$/caja 18 kilos"
$ws.Range("R159").Value = "Región de O'Higgins"
$ws.Range("S159").Value = 403
$ws.Range("T159").Value = 18

# New row 160: Terminal Hortofrutícola Agro Chillán, Ñuble, Mandarina - Clemenuless, Segunda
$ws.Range("A160").Value = 7
$ws.Range("B160").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C160").Value = "Ñuble"
$ws.Range("D160").Value = 44753
$ws.Range("E160").Value = 16
$ws.Range("F160").Value = "Fruta"
$ws.Range("G160").Value = 100102
$ws.Range("H160").Value = "Cítricos"
$ws.Range("I160").Value = 100102004
$ws.Range("J160").Value = "Mandarina"
$ws.Range("K160").Value = "Clemenuless"
$ws.Range("L160").Value = "Segunda"
$ws.Range("M160").Value = 120
$ws.Range("N160").Value = 6000
$ws.Range("O160").Value = 6500
$ws.Range("P160").Value = 6250
$ws.Range("Q160").Value = "$/caja 18 kilos"
$ws.Range("R160").Value = "Región de O'Higgins"
$ws.Range("S160").Value = 347
$ws.Range("T160").Value = 18
